$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.902115927284317
$ws.Range("C2").Value = -8.648538155816993
$ws.Range("D2").Value = -0.9391415561608464
$ws.Range("E2").Value = 0.6437711047260652
$ws.Range("F2").Value = -2.322515887623191
$ws.Range("G2").Value = -0.05848473295872768
$ws.Range("H2").Value = -0.5096388878854663
$ws.Range("I2").Value = -0.6500029513994404
$ws.Range("J2").Value = -0.1237639405372229
$ws.Range("K2").Value = -0.2655928872764969
$ws.Range("B3").Value = -10.89834099542839
$ws.Range("C3").Value = -3.188944395772239
$ws.Range("D3").Value = -1.606031734885327
$ws.Range("E3").Value = -4.572318727234583
$ws.Range("F3").Value = -2.30828757257012
$ws.Range("G3").Value = -2.759441727496859
$ws.Range("H3").Value = -2.899805791010833
$ws.Range("I3").Value = -2.373566780148615
$ws.Range("J3").Value = -2.515395726887889
$ws.Range("K3").Value = -1.872421105485825
$ws.Range("B4").Value = -10.65809472662953
$ws.Range("C4").Value = -9.07518206574262
$ws.Range("D4").Value = -12.04146905809188
$ws.Range("E4").Value = -9.777437903427414
$ws.Range("F4").Value = -10.22859205835415
$ws.Range("G4").Value = -10.36895612186813
$ws.Range("H4").Value = -9.842717111005909
$ws.Range("I4").Value = -9.984546057745183
$ws.Range("J4").Value = -9.341571436343118
$ws.Range("K4").Value = -9.217158730013535
$ws.Range("B5").Value = 9.096079734394849
$ws.Range("C5").Value = 6.129792742045593
$ws.Range("D5").Value = 8.393823896710055
$ws.Range("E5").Value = 7.942669741783317
$ws.Range("F5").Value = 7.802305678269343
$ws.Range("G5").Value = 8.32854468913156
$ws.Range("H5").Value = 8.186715742392286
$ws.Range("I5").Value = 8.82969036379435
$ws.Range("J5").Value = 8.954103070123933
$ws.Range("K5").Value = 8.556472864086315
$ws.Range("B6").Value = -2.009878804933557
$ws.Range("C6").Value = 0.2541523497309064
$ws.Range("D6").Value = -0.1970018051958322
$ws.Range("E6").Value = -0.3373658687098063
$ws.Range("F6").Value = 0.1888731421524112
$ws.Range("G6").Value = 0.0470441954131372
$ws.Range("H6").Value = 0.6900188168152011
$ws.Range("I6").Value = 0.8144315231447835
$ws.Range("J6").Value = 0.4168013171071663
$ws.Range("K6").Value = 0.3406326034063205
$ws.Range("B7").Value = -1.893418122068486
$ws.Range("C7").Value = -2.344572276995224
$ws.Range("D7").Value = -2.484936340509199
$ws.Range("E7").Value = -1.958697329646981
$ws.Range("F7").Value = -2.100526276386255
$ws.Range("G7").Value = -1.457551654984191
$ws.Range("H7").Value = -1.333138948654609
$ws.Range("I7").Value = -1.730769154692226
$ws.Range("J7").Value = -1.806937868393072
$ws.Range("K7").Value = -1.53061581027015
$ws.Range("B8").Value = 1.095457709528105
$ws.Range("C8").Value = 0.9550936460141312
$ws.Range("D8").Value = 1.481332656876349
$ws.Range("E8").Value = 1.339503710137075
$ws.Range("F8").Value = 1.982478331539139
$ws.Range("G8").Value = 2.106891037868721
$ws.Range("H8").Value = 1.709260831831104
$ws.Range("I8").Value = 1.633092118130258
$ws.Range("J8").Value = 1.90941417625318
$ws.Range("K8").Value = 1.738194353941456
$ws.Range("B9").Value = 1.016267824428332
$ws.Range("C9").Value = 1.54250683529055
$ws.Range("D9").Value = 1.400677888551276
$ws.Range("E9").Value = 2.04365250995334
$ws.Range("F9").Value = 2.168065216282922
$ws.Range("G9").Value = 1.770435010245305
$ws.Range("H9").Value = 1.694266296544459
$ws.Range("I9").Value = 1.970588354667381
$ws.Range("J9").Value = 1.799368532355657
$ws.Range("K9").Value = 1.952515937864398
$ws.Range("B10").Value = -0.4989491018657047
$ws.Range("C10").Value = -0.6407780486049788
$ws.Range("D10").Value = 0.002196572797085183
$ws.Range("E10").Value = 0.1266092791266676
$ws.Range("F10").Value = -0.2710209269109496
$ws.Range("G10").Value = -0.3471896406117954
$ws.Range("H10").Value = -0.07086758248887381
$ws.Range("I10").Value = -0.2420874048005978
$ws.Range("J10").Value = -0.08893999929185659
$ws.Range("K10").Value = -0.3663687737149753
$ws.Range("B11").Value = -0.05818351157133772
$ws.Range("C11").Value = 0.5847911098307261
$ws.Range("D11").Value = 0.7092038161603086
$ws.Range("E11").Value = 0.3115736101226914
$ws.Range("F11").Value = 0.2354048964218456
$ws.Range("G11").Value = 0.5117269545447671
$ws.Range("H11").Value = 0.3405071322330432
$ws.Range("I11").Value = 0.4936545377417844
$ws.Range("J11").Value = 0.2162257633186657
$ws.Range("K11").Value = 0.05323806995971928
$ws.Range("B12").Value = 0.4891160690213684
$ws.Range("C12").Value = 0.6135287753509509
$ws.Range("D12").Value = 0.2158985693133336
$ws.Range("E12").Value = 0.1397298556124878
$ws.Range("F12").Value = 0.4160519137354094
$ws.Range("G12").Value = 0.2448320914236854
$ws.Range("H12").Value = 0.3979794969324266
$ws.Range("I12").Value = 0.1205507225093079
$ws.Range("J12").Value = -0.04243697084963852
$ws.Range("K12").Value = -0.1665195462441563
$ws.Range("B13").Value = 0.8739478123496736
$ws.Range("C13").Value = 0.4763176063120564
$ws.Range("D13").Value = 0.4001488926112106
$ws.Range("E13").Value = 0.6764709507341322
$ws.Range("F13").Value = 0.5052511284224082
$ws.Range("G13").Value = 0.6583985339311494
$ws.Range("H13").Value = 0.3809697595080307
$ws.Range("I13").Value = 0.2179820661490843
$ws.Range("J13").Value = 0.09389949075456649
$ws.Range("K13").Value = 0.6100525277605273
$ws.Range("B14").Value = -0.3591373898464348
$ws.Range("C14").Value = -0.4353061035472806
$ws.Range("D14").Value = -0.158984045424359
$ws.Range("E14").Value = -0.330203867736083
$ws.Range("F14").Value = -0.1770564622273418
$ws.Range("G14").Value = -0.4544852366504605
$ws.Range("H14").Value = -0.617472930009407
$ws.Range("I14").Value = -0.7415555054039247
$ws.Range("J14").Value = -0.2254024683979639
$ws.Range("K14").Value = -0.4353061035472806
$ws.Range("B15").Value = -0.3352267436446591
$ws.Range("C15").Value = -0.0589046855217375
$ws.Range("D15").Value = -0.2301245078334615
$ws.Range("E15").Value = -0.07697710232472027
$ws.Range("F15").Value = -0.354405876747839
$ws.Range("G15").Value = -0.5173935701067854
$ws.Range("H15").Value = -0.6414761455013032
$ws.Range("I15").Value = -0.1253231084953424
$ws.Range("J15").Value = -0.3352267436446591
$ws.Range("K15").ClearContents()
$ws.Range("B16").Value = 0.2948818205579588
$ws.Range("C16").Value = 0.1236619982462347
$ws.Range("D16").Value = 0.276809403754976
$ws.Range("E16").Value = -0.0006193706681427817
$ws.Range("F16").Value = -0.1636070640270892
$ws.Range("G16").Value = -0.287689639421607
$ws.Range("H16").Value = 0.2284633975843539
$ws.Range("I16").Value = 0.01855976243503714
$ws.Range("J16").ClearContents()
$ws.Range("B17").Value = -0.0245153921862106
$ws.Range("C17").Value = 0.1286320133225306
$ws.Range("D17").Value = -0.1487967611005881
$ws.Range("E17").Value = -0.3117844544595345
$ws.Range("F17").Value = -0.4358670298540523
$ws.Range("G17").Value = 0.08028600715190851
$ws.Range("H17").Value = -0.1296176279974082
$ws.Range("I17").ClearContents()
$ws.Range("B18").Value = -0.02881397568162436
$ws.Range("C18").Value = -0.3062427501047431
$ws.Range("D18").Value = -0.4692304434636895
$ws.Range("E18").Value = -0.5933130188582073
$ws.Range("F18").Value = -0.07715998185224648
$ws.Range("G18").Value = -0.2870636170015632
$ws.Range("H18").ClearContents()
$ws.Range("B19").Value = 0.1944167064213277
$ws.Range("C19").Value = 0.0314290130623813
$ws.Range("D19").Value = -0.09265356233213651
$ws.Range("E19").Value = 0.4234994746738243
$ws.Range("F19").Value = 0.2135958395245076
$ws.Range("G19").ClearContents()
$ws.Range("B20").Value = -0.2489288674730878
$ws.Range("C20").Value = -0.3730114428676057
$ws.Range("D20").Value = 0.1431415941383551
$ws.Range("E20").Value = -0.06676204101096155
$ws.Range("F20").ClearContents()
$ws.Range("B21").Value = -0.201036585022594
$ws.Range("C21").Value = 0.3151164519833668
$ws.Range("D21").Value = 0.1052128168340501
$ws.Range("E21").ClearContents()
$ws.Range("B22").Value = 0.009253912237035311
$ws.Range("C22").Value = -0.2006497229122814
$ws.Range("D22").ClearContents()
$ws.Range("B23").Value = 0.4116802297750048
$ws.Range("C23").ClearContents()
$ws.Range("B24").ClearContents()
